# Apply updated dSF (column F) values onto Sheet1.
# These reflect a "repull data, push all data, mean calculation" refresh
# where only column F (dSF) values changed for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2  = -3
    3  = -3
    4  = 0
    5  = -1
    6  = 6
    7  = -5
    9  = 1
    15 = 0
    16 = -5
    18 = 1
    19 = -6
    21 = -10
    22 = -3
    25 = -2
    26 = -7
    28 = 1
    33 = -2
    35 = -1
    37 = -2
    38 = -3
    40 = 7
    41 = 2
    42 = 6
    43 = 4
    44 = -7
    46 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
